$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 134.93333
$ws.Range("I11").Value = 134.93333
$ws.Range("K11").Value = 134.93333
$ws.Range("M11").Value = 5.066669999999988
$ws.Range("H17").Value = 2442.6943
$ws.Range("J17").Value = 2442.6943
$ws.Range("L17").Value = 7328.0829
$ws.Range("N17").Value = -7664.0829
$ws.Range("H32").Value = 9491.166999999999
$ws.Range("I32").Value = 7892.3335
$ws.Range("J32").Value = 10024.111
$ws.Range("K32").Value = 7892.3335
$ws.Range("L32").Value = 10024.111
$ws.Range("M32").Value = -7566.3335
$ws.Range("N32").Value = -10676.111
$ws.Range("H37").Value = 999
$ws.Range("I37").Value = 999
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2997
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2871
$ws.Range("N37").ClearContents()
$ws.Range("H40").Value = 5583.067
$ws.Range("I40").Value = 4299.375
$ws.Range("J40").Value = 7050.143
$ws.Range("K40").Value = 4299.375
$ws.Range("L40").Value = 7050.143
$ws.Range("M40").Value = -4124.375
$ws.Range("N40").Value = -7400.143
$ws.Range("H41").Value = 2055.4443
$ws.Range("I41").Value = 583.1667
$ws.Range("K41").Value = 583.1667
$ws.Range("M41").Value = -143.1667
$ws.Range("H43").Value = 12766.923
$ws.Range("I43").Value = 10981.833
$ws.Range("K43").Value = 10981.833
$ws.Range("M43").Value = -10912.833
$ws.Range("H52").Value = 1625
$ws.Range("J52").Value = 700
$ws.Range("L52").Value = 2100
$ws.Range("N52").Value = -2420
$ws.Range("H55").Value = 240.08333
$ws.Range("I55").Value = 307.33334
$ws.Range("J55").Value = 38.333332
$ws.Range("K55").Value = 307.33334
$ws.Range("L55").Value = 38.333332
$ws.Range("M55").Value = -93.33334000000002
$ws.Range("N55").Value = -466.333332
$ws.Range("H64").Value = 5760.8
$ws.Range("I64").Value = 5100.25
$ws.Range("K64").Value = 5100.25
$ws.Range("M64").Value = -4852.25
$ws.Range("H67").Value = 5760.8
$ws.Range("I67").Value = 5100.25
$ws.Range("K67").Value = 5100.25
$ws.Range("M67").Value = -4242.25
$ws.Range("H86").Value = 3297.889
$ws.Range("I86").Value = 1530.3334
$ws.Range("K86").Value = 1530.3334
$ws.Range("M86").Value = -407.3334
$ws.Range("H89").Value = 3297.889
$ws.Range("I89").Value = 1530.3334
$ws.Range("K89").Value = 7651.666999999999
$ws.Range("M89").Value = -2035.666999999999
$ws.Range("H99").Value = 383.5
$ws.Range("J99").Value = 780
$ws.Range("L99").Value = 2340
$ws.Range("N99").Value = -5336
$ws.Range("H112").Value = 2026.2333
$ws.Range("J112").Value = 2026.2333
$ws.Range("L112").Value = 6078.699900000001
$ws.Range("N112").Value = -8294.6999
$ws.Range("H118").Value = 339248.34
$ws.Range("I118").Value = 339248.34
$ws.Range("K118").Value = 1017745.02
$ws.Range("M118").Value = -1016088.02
$ws.Range("H135").Value = 654.6667
$ws.Range("I135").Value = 589.95654
$ws.Range("J135").Value = 1026.75
$ws.Range("K135").Value = 5309.60886
$ws.Range("L135").Value = 9240.75
$ws.Range("M135").Value = -2774.60886
$ws.Range("N135").Value = -14310.75
$ws.Range("H137").Value = 10667.5
$ws.Range("I137").Value = 1797.4286
$ws.Range("J137").Value = 21624.646
$ws.Range("K137").Value = 5392.2858
$ws.Range("L137").Value = 64873.938
$ws.Range("M137").Value = -2842.2858
$ws.Range("N137").Value = -69973.93799999999
$ws.Range("H138").Value = 2579.831
$ws.Range("I138").Value = 1542.8298
$ws.Range("J138").Value = 4204.467
$ws.Range("K138").Value = 4628.4894
$ws.Range("L138").Value = 12613.401
$ws.Range("M138").Value = 511.5105999999996
$ws.Range("N138").Value = -22893.401
$ws.Range("H141").Value = 1739.875
$ws.Range("I141").Value = 1769.6666
$ws.Range("J141").Value = 1650.5
$ws.Range("K141").Value = 5308.9998
$ws.Range("L141").Value = 4951.5
$ws.Range("M141").Value = -128.9997999999996
$ws.Range("N141").Value = -15311.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5414492
$ws.Range("I2").Value = 6573812
$ws.Range("J2").Value = 4333
$ws.Range("K2").Value = 6573812
$ws.Range("L2").Value = 4333
$ws.Range("M2").Value = -6573699
$ws.Range("N2").Value = -4559
$ws.Range("H5").Value = 234.83333
$ws.Range("J5").Value = 13.75
$ws.Range("L5").Value = 13.75
$ws.Range("N5").Value = -237.75
$ws.Range("H32").Value = 2960061.2
$ws.Range("J32").Value = 9799.6
$ws.Range("L32").Value = 9799.6
$ws.Range("N32").Value = -10373.6
$ws.Range("H45").Value = 36666.668
$ws.Range("I45").Value = 52000
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 52000
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -51623
$ws.Range("N45").Value = -6754
$ws.Range("H61").Value = 6686.8076
$ws.Range("I61").Value = 6264.1396
$ws.Range("J61").Value = 8706.223
$ws.Range("K61").Value = 6264.1396
$ws.Range("L61").Value = 8706.223
$ws.Range("M61").Value = -6052.1396
$ws.Range("N61").Value = -9130.223
$ws.Range("H63").Value = 6913.794
$ws.Range("I63").Value = 743
$ws.Range("J63").Value = 9865.044
$ws.Range("K63").Value = 743
$ws.Range("L63").Value = 9865.044
$ws.Range("M63").Value = -57
$ws.Range("N63").Value = -11237.044
$ws.Range("H66").Value = 6913.794
$ws.Range("I66").Value = 743
$ws.Range("J66").Value = 9865.044
$ws.Range("K66").Value = 3715
$ws.Range("L66").Value = 49325.22
$ws.Range("M66").Value = -283
$ws.Range("N66").Value = -56189.22
$ws.Range("H74").Value = 591760.9
$ws.Range("I74").Value = 834865.8
$ws.Range("K74").Value = 834865.8
$ws.Range("M74").Value = -833991.8
$ws.Range("H77").Value = 591760.9
$ws.Range("I77").Value = 834865.8
$ws.Range("K77").Value = 4174329
$ws.Range("M77").Value = -4169961
$ws.Range("H95").Value = 17375
$ws.Range("J95").Value = 17166.666
$ws.Range("L95").Value = 17166.666
$ws.Range("N95").Value = -22658.666
$ws.Range("H116").Value = 5414492
$ws.Range("I116").Value = 6573812
$ws.Range("J116").Value = 4333
$ws.Range("K116").Value = 6573812
$ws.Range("L116").Value = 4333
$ws.Range("M116").Value = -6571518
$ws.Range("N116").Value = -8921
$ws.Range("H132").Value = 4392003
$ws.Range("I132").Value = 6177007
$ws.Range("J132").Value = 10628.728
$ws.Range("K132").Value = 18531021
$ws.Range("L132").Value = 31886.184
$ws.Range("M132").Value = -18528491
$ws.Range("N132").Value = -36946.18399999999
$ws.Range("H135").Value = 137333
$ws.Range("J135").Value = 137333
$ws.Range("L135").Value = 137333
$ws.Range("N135").Value = -147473
$ws.Range("H136").Value = 6686.8076
$ws.Range("I136").Value = 6264.1396
$ws.Range("J136").Value = 8706.223
$ws.Range("K136").Value = 18792.4188
$ws.Range("L136").Value = 26118.669
$ws.Range("M136").Value = -16242.4188
$ws.Range("N136").Value = -31218.669
$ws.Range("H139").Value = 67618.836
$ws.Range("J139").Value = 67142.60000000001
$ws.Range("L139").Value = 67142.60000000001
$ws.Range("N139").Value = -77422.60000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5414492
$ws.Range("I3").Value = 6573812
$ws.Range("J3").Value = 4333
$ws.Range("K3").Value = 6573812
$ws.Range("L3").Value = 4333
$ws.Range("M3").Value = -6573698
$ws.Range("N3").Value = -4561
$ws.Range("H4").Value = 234.83333
$ws.Range("J4").Value = 13.75
$ws.Range("L4").Value = 13.75
$ws.Range("N4").Value = -243.75
$ws.Range("H24").Value = 4875
$ws.Range("I24").Value = 4833.3335
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 4833.3335
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -4598.3335
$ws.Range("N24").Value = -5470
$ws.Range("H53").Value = 70000
$ws.Range("J53").Value = 70000
$ws.Range("L53").Value = 70000
$ws.Range("N53").Value = -71148
$ws.Range("H64").Value = 1171.2858
$ws.Range("I64").Value = 611.1667
$ws.Range("K64").Value = 611.1667
$ws.Range("M64").Value = -386.1667
$ws.Range("H67").Value = 1171.2858
$ws.Range("I67").Value = 611.1667
$ws.Range("K67").Value = 611.1667
$ws.Range("M67").Value = 168.8333
$ws.Range("H86").Value = 287501
$ws.Range("I86").Value = 2001.6
$ws.Range("K86").Value = 2001.6
$ws.Range("M86").Value = -878.5999999999999
$ws.Range("H89").Value = 287501
$ws.Range("I89").Value = 2001.6
$ws.Range("K89").Value = 10008
$ws.Range("M89").Value = -4392
$ws.Range("H102").Value = 29185.334
$ws.Range("I102").Value = 29185.334
$ws.Range("K102").Value = 29185.334
$ws.Range("M102").Value = -25940.334
$ws.Range("H107").Value = 1417.7778
$ws.Range("I107").Value = 1238.8334
$ws.Range("J107").Value = 1775.6666
$ws.Range("K107").Value = 1238.8334
$ws.Range("L107").Value = 1775.6666
$ws.Range("M107").Value = 681.1666
$ws.Range("N107").Value = -5615.6666
$ws.Range("H134").Value = 6493.375
$ws.Range("I134").Value = 4935.4443
$ws.Range("J134").Value = 8496.429
$ws.Range("K134").Value = 14806.3329
$ws.Range("L134").Value = 25489.287
$ws.Range("M134").Value = -12271.3329
$ws.Range("N134").Value = -30559.287

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 200.4
$ws.Range("I7").Value = 243.81818
$ws.Range("J7").Value = 81
$ws.Range("K7").Value = 243.81818
$ws.Range("L7").Value = 81
$ws.Range("M7").Value = -130.81818
$ws.Range("N7").Value = -307
$ws.Range("H22").Value = 2795.25
$ws.Range("I22").Value = 749.6667
$ws.Range("J22").Value = 4022.6
$ws.Range("K22").Value = 749.6667
$ws.Range("L22").Value = 4022.6
$ws.Range("M22").Value = -399.6667
$ws.Range("N22").Value = -4722.6
$ws.Range("H31").Value = 29415644
$ws.Range("I31").Value = 50002132
$ws.Range("J31").Value = 6373.5
$ws.Range("K31").Value = 50002132
$ws.Range("L31").Value = 6373.5
$ws.Range("M31").Value = -50001837
$ws.Range("N31").Value = -6963.5
$ws.Range("H34").Value = 29415644
$ws.Range("I34").Value = 50002132
$ws.Range("J34").Value = 6373.5
$ws.Range("K34").Value = 50002132
$ws.Range("L34").Value = 6373.5
$ws.Range("M34").Value = -50001930
$ws.Range("N34").Value = -6777.5
$ws.Range("H58").Value = 10011.421
$ws.Range("I58").Value = 14542.875
$ws.Range("J58").Value = 6715.8184
$ws.Range("K58").Value = 14542.875
$ws.Range("L58").Value = 6715.8184
$ws.Range("M58").Value = -14339.875
$ws.Range("N58").Value = -7121.8184
$ws.Range("H86").Value = 8292
$ws.Range("I86").Value = 8099.75
$ws.Range("J86").Value = 8599.6
$ws.Range("K86").Value = 8099.75
$ws.Range("L86").Value = 8599.6
$ws.Range("M86").Value = -6976.75
$ws.Range("N86").Value = -10845.6
$ws.Range("H89").Value = 8292
$ws.Range("I89").Value = 8099.75
$ws.Range("J89").Value = 8599.6
$ws.Range("K89").Value = 40498.75
$ws.Range("L89").Value = 42998
$ws.Range("M89").Value = -34882.75
$ws.Range("N89").Value = -54230
$ws.Range("H99").Value = 9308.333000000001
$ws.Range("I99").Value = 10642.471
$ws.Range("K99").Value = 10642.471
$ws.Range("M99").Value = -9144.471
$ws.Range("H107").Value = 1133.4
$ws.Range("I107").Value = 686.8333
$ws.Range("K107").Value = 686.8333
$ws.Range("M107").Value = 1233.1667
$ws.Range("H112").Value = 74092.60000000001
$ws.Range("J112").Value = 74092.60000000001
$ws.Range("L112").Value = 74092.60000000001
$ws.Range("N112").Value = -77046.60000000001
$ws.Range("H122").Value = 68036.2
$ws.Range("I122").Value = 125921
$ws.Range("J122").Value = 1882.1428
$ws.Range("K122").Value = 377763
$ws.Range("L122").Value = 5646.428400000001
$ws.Range("M122").Value = -375313
$ws.Range("N122").Value = -10546.4284
$ws.Range("H126").Value = 9308.333000000001
$ws.Range("I126").Value = 10642.471
$ws.Range("K126").Value = 31927.413
$ws.Range("M126").Value = -29457.413
$ws.Range("H132").Value = 4583.4136
$ws.Range("I132").Value = 3876.5789
$ws.Range("K132").Value = 11629.7367
$ws.Range("M132").Value = -9099.736699999999
$ws.Range("H133").Value = 52250
$ws.Range("J133").Value = 52250
$ws.Range("L133").Value = 52250
$ws.Range("N133").Value = -57310
$ws.Range("H134").Value = 4871
$ws.Range("I134").Value = 4458
$ws.Range("K134").Value = 13374
$ws.Range("M134").Value = -10839
$ws.Range("H135").Value = 105275.29
$ws.Range("J135").Value = 105275.29
$ws.Range("L135").Value = 105275.29
$ws.Range("N135").Value = -115415.29
$ws.Range("H136").Value = 10011.421
$ws.Range("I136").Value = 14542.875
$ws.Range("J136").Value = 6715.8184
$ws.Range("K136").Value = 43628.625
$ws.Range("L136").Value = 20147.4552
$ws.Range("M136").Value = -41078.625
$ws.Range("N136").Value = -25247.4552

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2239288.2
$ws.Range("I4").Value = 3400173.8
$ws.Range("J4").Value = 72302.2
$ws.Range("K4").Value = 10200521.4
$ws.Range("L4").Value = 216906.6
$ws.Range("M4").Value = -10200409.4
$ws.Range("N4").Value = -217130.6
$ws.Range("H12").Value = 287.33334
$ws.Range("J12").Value = 331.23077
$ws.Range("L12").Value = 993.69231
$ws.Range("N12").Value = -1339.69231
$ws.Range("H33").Value = 22.866667
$ws.Range("I33").Value = 25
$ws.Range("J33").Value = 18.6
$ws.Range("K33").Value = 150
$ws.Range("L33").Value = 111.6
$ws.Range("M33").Value = 133
$ws.Range("N33").Value = -677.6
$ws.Range("H44").Value = 2698.1667
$ws.Range("I44").Value = 396.33334
$ws.Range("J44").Value = 5000
$ws.Range("K44").Value = 1189.00002
$ws.Range("L44").Value = 15000
$ws.Range("M44").Value = -791.0000199999999
$ws.Range("N44").Value = -15796
$ws.Range("H63").Value = 7866.6665
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 15000
$ws.Range("M63").Value = -14251
$ws.Range("H66").Value = 7866.6665
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 45000
$ws.Range("M66").Value = -41256
$ws.Range("H81").Value = 963
$ws.Range("I81").Value = 963
$ws.Range("K81").Value = 2889
$ws.Range("M81").Value = -1766
$ws.Range("H84").Value = 963
$ws.Range("I84").Value = 963
$ws.Range("K84").Value = 8667
$ws.Range("M84").Value = -3051
$ws.Range("H98").Value = 245.6
$ws.Range("I98").Value = 199.5
$ws.Range("J98").Value = 276.33334
$ws.Range("K98").Value = 598.5
$ws.Range("L98").Value = 829.0000200000001
$ws.Range("M98").Value = 899.5
$ws.Range("N98").Value = -3825.00002
$ws.Range("H139").Value = 3462
$ws.Range("I139").Value = 3144.4285
$ws.Range("K139").Value = 9433.2855
$ws.Range("M139").Value = -4293.2855
$ws.Range("H140").Value = 22728768
$ws.Range("I140").Value = 41667570
$ws.Range("J140").Value = 2206.6
$ws.Range("K140").Value = 125002710
$ws.Range("L140").Value = 6619.799999999999
$ws.Range("M140").Value = -124997530
$ws.Range("N140").Value = -16979.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 7288.6665
$ws.Range("J19").Value = 7449.75
$ws.Range("L19").Value = 7449.75
$ws.Range("N19").Value = -8025.75
$ws.Range("H41").Value = 654.5
$ws.Range("I41").Value = 654.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 654.5
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("M41").Value = -299.5
$ws.Range("H80").Value = 3532.7144
$ws.Range("I80").Value = 2947.75
$ws.Range("J80").Value = 4312.6665
$ws.Range("K80").Value = 2947.75
$ws.Range("L80").Value = 4312.6665
$ws.Range("M80").Value = -1949.75
$ws.Range("N80").Value = -6308.6665
$ws.Range("H83").Value = 3532.7144
$ws.Range("I83").Value = 2947.75
$ws.Range("J83").Value = 4312.6665
$ws.Range("K83").Value = 14738.75
$ws.Range("L83").Value = 21563.3325
$ws.Range("M83").Value = -9746.75
$ws.Range("N83").Value = -31547.3325
$ws.Range("H107").Value = 1513.3846
$ws.Range("I107").Value = 1577.3334
$ws.Range("J107").Value = 1458.5714
$ws.Range("K107").Value = 1577.3334
$ws.Range("L107").Value = 1458.5714
$ws.Range("M107").Value = 342.6666
$ws.Range("N107").Value = -5298.5714
$ws.Range("H132").Value = 7303.421
$ws.Range("I132").Value = 4256.5
$ws.Range("J132").Value = 10688.889
$ws.Range("K132").Value = 12769.5
$ws.Range("L132").Value = 32066.667
$ws.Range("M132").Value = -10239.5
$ws.Range("N132").Value = -37126.667
$ws.Range("H133").Value = 210000
$ws.Range("J133").Value = 210000
$ws.Range("L133").Value = 210000
$ws.Range("N133").Value = -220120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 378998.5
$ws.Range("I7").Value = 378998.5
$ws.Range("K7").Value = 378998.5
$ws.Range("M7").Value = -378886.5
$ws.Range("H16").Value = 1006.8571
$ws.Range("I16").Value = 1006.8571
$ws.Range("K16").Value = 1006.8571
$ws.Range("M16").Value = -836.8570999999999
$ws.Range("H17").Value = 8000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 8000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -8340
$ws.Range("H22").Value = 3159
$ws.Range("J22").Value = 3503.3157
$ws.Range("L22").Value = 3503.3157
$ws.Range("N22").Value = -4093.3157
$ws.Range("H27").Value = 3159
$ws.Range("J27").Value = 3503.3157
$ws.Range("L27").Value = 3503.3157
$ws.Range("N27").Value = -3717.3157
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H40").Value = 18332.834
$ws.Range("I40").Value = 16199.6
$ws.Range("K40").Value = 16199.6
$ws.Range("M40").Value = -16063.6
$ws.Range("H46").Value = 4323.05
$ws.Range("I46").Value = 1149.2858
$ws.Range("K46").Value = 1149.2858
$ws.Range("M46").Value = -961.2858000000001
$ws.Range("H61").Value = 4706.5557
$ws.Range("J61").Value = 7313.75
$ws.Range("L61").Value = 7313.75
$ws.Range("N61").Value = -7717.75
$ws.Range("H68").Value = 3266.7778
$ws.Range("I68").Value = 2709.3635
$ws.Range("K68").Value = 2709.3635
$ws.Range("M68").Value = -1960.3635
$ws.Range("H71").Value = 3266.7778
$ws.Range("I71").Value = 2709.3635
$ws.Range("K71").Value = 13546.8175
$ws.Range("M71").Value = -9802.817499999999
$ws.Range("H82").Value = 2076.4783
$ws.Range("I82").Value = 1968.625
$ws.Range("K82").Value = 1968.625
$ws.Range("M82").Value = -1607.625
$ws.Range("H85").Value = 2076.4783
$ws.Range("I85").Value = 1968.625
$ws.Range("K85").Value = 1968.625
$ws.Range("M85").Value = -720.625
$ws.Range("H93").Value = 1112.3636
$ws.Range("I93").Value = 1070.8889
$ws.Range("J93").Value = 1299
$ws.Range("K93").Value = 1070.8889
$ws.Range("L93").Value = 1299
$ws.Range("M93").Value = 177.1111000000001
$ws.Range("N93").Value = -3795
$ws.Range("H113").Value = 4706.5557
$ws.Range("J113").Value = 7313.75
$ws.Range("L113").Value = 7313.75
$ws.Range("N113").Value = -11653.75
$ws.Range("H122").Value = 100005496
$ws.Range("I122").Value = 125004870
$ws.Range("K122").Value = 375014610
$ws.Range("M122").Value = -375012160
$ws.Range("H126").Value = 378998.5
$ws.Range("I126").Value = 378998.5
$ws.Range("K126").Value = 1136995.5
$ws.Range("M126").Value = -1134525.5
$ws.Range("H132").Value = 4127.2173
$ws.Range("I132").Value = 2161.5386
$ws.Range("J132").Value = 6682.6
$ws.Range("K132").Value = 6484.6158
$ws.Range("L132").Value = 20047.8
$ws.Range("M132").Value = -3954.6158
$ws.Range("N132").Value = -25107.8
$ws.Range("H136").Value = 3856.1155
$ws.Range("I136").Value = 2845.6667
$ws.Range("J136").Value = 8100
$ws.Range("K136").Value = 8537.000100000001
$ws.Range("L136").Value = 24300
$ws.Range("M136").Value = -5987.000100000001
$ws.Range("N136").Value = -29400

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 204899.8
$ws.Range("I2").Value = 5750
$ws.Range("J2").Value = 337666.34
$ws.Range("K2").Value = 5750
$ws.Range("L2").Value = 337666.34
$ws.Range("M2").Value = -5638
$ws.Range("N2").Value = -337890.34
$ws.Range("H4").Value = 642919.5
$ws.Range("I4").Value = 69833.336
$ws.Range("J4").Value = 986771.2
$ws.Range("K4").Value = 69833.336
$ws.Range("L4").Value = 986771.2
$ws.Range("M4").Value = -69720.336
$ws.Range("N4").Value = -986997.2
$ws.Range("H40").Value = 59999.5
$ws.Range("I40").Value = 59999
$ws.Range("K40").Value = 59999
$ws.Range("M40").Value = -59850
$ws.Range("H62").Value = 9161
$ws.Range("I62").Value = 8861.888999999999
$ws.Range("K62").Value = 8861.888999999999
$ws.Range("M62").Value = -8237.888999999999
$ws.Range("H65").Value = 9161
$ws.Range("I65").Value = 8861.888999999999
$ws.Range("K65").Value = 44309.44499999999
$ws.Range("M65").Value = -41189.44499999999
$ws.Range("H70").Value = 35666
$ws.Range("I70").Value = 33499
$ws.Range("J70").Value = 40000
$ws.Range("K70").Value = 33499
$ws.Range("L70").Value = 40000
$ws.Range("M70").Value = -33184
$ws.Range("N70").Value = -40630
$ws.Range("H73").Value = 35666
$ws.Range("I73").Value = 33499
$ws.Range("J73").Value = 40000
$ws.Range("K73").Value = 33499
$ws.Range("L73").Value = 40000
$ws.Range("M73").Value = -32407
$ws.Range("N73").Value = -42184
$ws.Range("H81").Value = 12478.728
$ws.Range("I81").Value = 3664.2144
$ws.Range("J81").Value = 15488.561
$ws.Range("K81").Value = 7328.4288
$ws.Range("L81").Value = 30977.122
$ws.Range("M81").Value = -6267.4288
$ws.Range("N81").Value = -33099.122
$ws.Range("H84").Value = 12478.728
$ws.Range("I84").Value = 3664.2144
$ws.Range("J84").Value = 15488.561
$ws.Range("K84").Value = 36642.144
$ws.Range("L84").Value = 154885.61
$ws.Range("M84").Value = -31338.144
$ws.Range("N84").Value = -165493.61
$ws.Range("H96").Value = 4176.2666
$ws.Range("I96").Value = 2805.7144
$ws.Range("K96").Value = 2805.7144
$ws.Range("M96").Value = -1432.7144
$ws.Range("H107").Value = 4380.2104
$ws.Range("I107").Value = 4013.2942
$ws.Range("J107").Value = 7499
$ws.Range("K107").Value = 12039.8826
$ws.Range("L107").Value = 22497
$ws.Range("M107").Value = -10119.8826
$ws.Range("N107").Value = -26337
$ws.Range("H113").Value = 925.0357
$ws.Range("I113").Value = 736.85
$ws.Range("J113").Value = 1395.5
$ws.Range("K113").Value = 2210.55
$ws.Range("L113").Value = 4186.5
$ws.Range("M113").Value = -40.55000000000018
$ws.Range("N113").Value = -8526.5
$ws.Range("H122").Value = 2006.9474
$ws.Range("I122").Value = 1651.75
$ws.Range("J122").Value = 3901.3333
$ws.Range("K122").Value = 4955.25
$ws.Range("L122").Value = 11703.9999
$ws.Range("M122").Value = -2505.25
$ws.Range("N122").Value = -16603.9999
$ws.Range("H132").Value = 10872903
$ws.Range("I132").Value = 15153981
$ws.Range("J132").Value = 5550.615
$ws.Range("K132").Value = 45461943
$ws.Range("L132").Value = 16651.845
$ws.Range("M132").Value = -45459413
$ws.Range("N132").Value = -21711.845
$ws.Range("H136").Value = 3675.8
$ws.Range("I136").Value = 1055.2307
$ws.Range("J136").Value = 8542.571
$ws.Range("K136").Value = 3165.6921
$ws.Range("L136").Value = 25627.713
$ws.Range("M136").Value = -615.6921000000002
$ws.Range("N136").Value = -30727.713

